# Release-Notes.xlsx update
# - "Folder Inventory": a new folder entry was captured, so it is inserted
#   as the new most-recent row (row 2), pushing the existing rows down by
#   one. The previous occurrence of that same folder further down the list
#   (it had not changed since its own last update) is removed so it is not
#   duplicated.
# - "Metadata": the generation timestamp and workflow run counter advance.
# - "Summary": the most-recent-update metric advances to match the new
#   top row of the Folder Inventory sheet.

$wb = $excel.ActiveWorkbook

$wsInventory = $wb.Worksheets.Item("Folder Inventory")
$wsMetadata  = $wb.Worksheets.Item("Metadata")
$wsSummary   = $wb.Worksheets.Item("Summary")

# --- Folder Inventory ---------------------------------------------------
# Insert a fresh row directly under the header so the new entry becomes
# the most-recent item, shifting every row below it down by one. Clear the
# formats the insert copies down from the header row so the new row keeps
# the plain (unstyled) look of the rest of the data rows.
$wsInventory.Range("A2:E2").Insert(-4121)
$wsInventory.Range("A2:E2").ClearFormats()

$wsInventory.Range("A2").Value = "Automate document processing by using Azure AI & OpenAI"
$wsInventory.Range("B2").Value = "Automate document processing by using Azure AI & OpenAI"
$wsInventory.Range("C2").Value = "2025-06-16 19:39:25 +0530"
$wsInventory.Range("D2").Value = 1
$wsInventory.Range("E2").Value = "Root"

# The same folder's earlier row (now at row 29 after the insert, formerly
# row 28) is the stale duplicate of the entry we just moved to the top --
# remove it so the folder only appears once.
$wsInventory.Rows.Item(29).Delete()

# --- Metadata -------------------------------------------------------------
$wsMetadata.Range("B3").Value = "2025-06-16 14:09:47 UTC"
# "Workflow Run" is stored as text ("10"/"11"), not a number -- a leading
# apostrophe keeps the COM layer from auto-coercing it to numeric 11, and
# resetting the style afterwards drops the quote-prefix formatting so the
# cell matches the other plain (unstyled) data cells.
$wsMetadata.Range("B5").Value = "'11"
$wsMetadata.Range("B5").Style = "Normal"

# --- Summary ----------------------------------------------------------------
$wsSummary.Range("B5").Value = "2025-06-16 19:39:25 +0530"
